$d = $word.ActiveDocument

# Locate the "Requisitos" list entry for LOQ4086, which is the anchor right
# before the block of paragraphs that must be removed (the blank line, the
# "Ver no Jupiter Salvar em pdf Salvar em docx" line, another blank line,
# and the page-break paragraph that precedes the final two tail paragraphs).
$r = $d.Content
$found = $r.Find.Execute("LOQ4086: Operações Unitárias II (Requisito fraco)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $anchorEnd = $r.End

    # Resolve the paragraph index that contains the end of the found text.
    $count = $d.Paragraphs.Count
    $anchorIndex = -1
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $anchorEnd -and $anchorEnd -le $p.Range.End) {
            $anchorIndex = $i
            break
        }
    }

    if ($anchorIndex -gt 0) {
        # The four paragraphs immediately following the anchor are the ones
        # removed by the edit:
        #   1) blank "Normal" paragraph
        #   2) "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph
        #   3) blank "Normal" paragraph
        #   4) blank paragraph with pageBreakBefore + jc=left
        $startParagraph = $d.Paragraphs.Item($anchorIndex + 1)
        $endParagraph = $d.Paragraphs.Item($anchorIndex + 4)

        $deleteRange = $d.Range($startParagraph.Range.Start, $endParagraph.Range.End)
        $deleteRange.Delete()
    }
}
